# Residencial plan config: add "200" plan (promo) and clear the
# "Sin_TotalPlay_TV" rows (crécelo) per commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plans")

# Row 5: bump the Megas plan value from 20 to 200 (new "200" plan)
$ws.Range("D5").Value = 200

# Rows 6-9: clear out the old Residencial/Sin_TotalPlay_TV plan rows
$ws.Range("B6:D6").ClearContents()
$ws.Range("B7:D7").ClearContents()
$ws.Range("B8:D8").ClearContents()
$ws.Range("B9:D9").ClearContents()

# Move the active selection to D10 (matches the saved view state)
$ws.Range("D10").Select()
